# Applying translations for project views
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the Spanish "stillNoFiles" translation
$ws.Range("B23").Value = "Todavía no se han añadido documentos."

# Append new translation rows for the loading spinner / title list fragments
$ws.Range("A36").Value = "fragments.loadingSpinner.loading"
$ws.Range("B36").Value = "Actualizando las visualizaciones. Un momento…"
$ws.Range("C36").Value = "Loading visualizations. Please wait…"

$ws.Range("A37").Value = "fragments.titleList.title"
$ws.Range("B37").Value = "Lista de documentos"
$ws.Range("C37").Value = "Document list"

$ws.Range("A38").Value = "fragments.titleList.find"
$ws.Range("B38").Value = "Encuentra documentos"
$ws.Range("C38").Value = "Find documents"

$ws.Range("A39").Value = "fragments.titleList.noDocuments"
$ws.Range("B39").Value = "No hay documentos disponibles."
$ws.Range("C39").Value = "There are no available documents."

# Reflect the new selection left behind in the source file
$ws.Range("C40").Select()
